$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, which pushes the former rows 118-173
# down to 119-174 (keeping all of their data intact).
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new record's data.
$ws.Cells.Item(118, 1).Value2 = 5
$ws.Cells.Item(118, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(118, 3).Value2 = "Maule"
$ws.Cells.Item(118, 4).Value2 = 44704
$ws.Cells.Item(118, 5).Value2 = 7
$ws.Cells.Item(118, 6).Value2 = 100112017
$ws.Cells.Item(118, 7).Value2 = "Apio"
$ws.Cells.Item(118, 8).Value2 = "Americana (o)"
$ws.Cells.Item(118, 9).Value2 = "Primera"
$ws.Cells.Item(118, 10).Value2 = 800
$ws.Cells.Item(118, 11).Value2 = 6500
$ws.Cells.Item(118, 12).Value2 = 6500
$ws.Cells.Item(118, 13).Value2 = 6500
$ws.Cells.Item(118, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(118, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(118, 16).Value2 = 1083
$ws.Cells.Item(118, 17).Value2 = 6
$ws.Cells.Item(118, 18).Value2 = "Hortaliza"
